# "Squishes neighboring text runs into one if compatible."
#
# Cell A3 used to hold the plain shared string "dolor". It becomes a single
# rich-text cell made of two runs: "dolor " in the default cell font, and
# "osa" in a new blue (RGB 0070C0) Calibri font - i.e. two neighbouring runs
# squished into one shared-string entry.
#
# Also moves the sheet's active-cell selection from A10 to A4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Register the new font (sz 11, RGB 0070C0, Calibri/minor-scheme) in the
# --- workbook's style table via a throwaway named cell style, without
# --- touching any worksheet cell. This mirrors the <font> Excel appends to
# --- styles.xml's <fonts> collection when a run uses a colour never seen
# --- before in the workbook.
$scratchStyle = $wb.Styles.Add("__scratch_font__")
$scratchStyle.Font.Color = 12611584   # RGB(0, 112, 192) = 0070C0
$scratchStyle.Delete()

# --- Rewrite A3 as rich text: "dolor " (unchanged formatting) + "osa" (blue).
$cell = $ws.Range("A3")
$cell.Value = "dolor osa"
$blueRun = $cell.Characters(7, 3)
$blueRun.Font.Color = 12611584        # RGB(0, 112, 192) = 0070C0
$blueRun.Font.Name = "Calibri"
$blueRun.Font.Size = 11

# --- Move the sheet selection from A10 to A4.
$ws.Range("A4").Select()
